# Append new Q&A paragraphs (67-74) after the existing last paragraph (66),
# keeping the same run formatting (sz=40 / szCs=40) used throughout the doc.

$d = $word.ActiveDocument

$newQuestions = @(
    "67) what is H2 database and how will you set it up in your application?",
    "68)how will you configure the multiple databases in your application ?",
    "69) What are Spring profiles and their uses.",
    "70)how will you reduce the application startup time in spring boot application ?",
    "71)What is lombook dependency and how it’s helpful for developers ?",
    "72)@ControllerAdvice vs @RestControllerAdvice.",
    "73) Handling exceptions globally in SB.",
    "74) How do you handle @Valid exceptions in SB?"
)

foreach ($q in $newQuestions) {
    $last = $d.Paragraphs.Last
    $r = $last.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newLast = $d.Paragraphs.Last
    $newLast.Range.Text = $q
}

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
